$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Symmetric matrix block (B2:E5) updated with new covariance values.
$ws.Range("C2").Value = 0.0000001711173756007756
$ws.Range("D2").Value = 0.000000005049777682986309
$ws.Range("E2").Value = 0.000000003847902911326734

$ws.Range("B3").Value = 0.0000001711173756007756
$ws.Range("D3").Value = 0.000000007043876027046637
$ws.Range("E3").Value = 2.218573043980729

$ws.Range("B4").Value = 0.000000005049777682986309
$ws.Range("C4").Value = 0.000000007043876027046637
$ws.Range("E4").Value = 0.8517350933034392

$ws.Range("B5").Value = 0.000000003847902911326734
$ws.Range("C5").Value = 2.218573043980729
$ws.Range("D5").Value = 0.8517350933034392
